# Update the NATMI LR-pairs sheet with recomputed TPM values.
# Rows 2-5: existing pairs (ECs -> {ECs,FAPs,MuSCs,Resolving-Mac}) get refreshed numbers.
# Rows 6-9: new pairs (Resolving-Mac -> {ECs,FAPs,MuSCs,Resolving-Mac}) are appended.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lta"
$ws.Range("C2").Value = "Tnfrsf1b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4993536666666667
$ws.Range("H2").Value = 1.498061
$ws.Range("I2").Value = 0.9584860631692095
$ws.Range("J2").Value = 0.9584860631692095
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 16.05260533333333
$ws.Range("N2").Value = 48.157816
$ws.Range("O2").Value = 0.1752915379534001
$ws.Range("P2").Value = 0.1752915379534001
$ws.Range("Q2").Value = 8.01592733275289
$ws.Range("R2").Value = 72.143345994776
$ws.Range("S2").Value = 0.1680144961198305
$ws.Range("T2").Value = 0.1680144961198305

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lta"
$ws.Range("C3").Value = "Tnfrsf1b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4993536666666667
$ws.Range("H3").Value = 1.498061
$ws.Range("I3").Value = 0.9584860631692095
$ws.Range("J3").Value = 0.9584860631692095
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.224257666666666
$ws.Range("N3").Value = 21.672773
$ws.Range("O3").Value = 0.07888758308485012
$ws.Range("P3").Value = 0.07888758308485012
$ws.Range("Q3").Value = 3.607459554794778
$ws.Range("R3").Value = 32.467135993153
$ws.Range("S3").Value = 0.07561264894393192
$ws.Range("T3").Value = 0.07561264894393192

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lta"
$ws.Range("C4").Value = "Tnfrsf1b"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4993536666666667
$ws.Range("H4").Value = 1.498061
$ws.Range("I4").Value = 0.9584860631692095
$ws.Range("J4").Value = 0.9584860631692095
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.098187333333334
$ws.Range("N4").Value = 6.294562000000001
$ws.Range("O4").Value = 0.02291182502385553
$ws.Range("P4").Value = 0.02291182502385553
$ws.Range("Q4").Value = 1.047737538253556
$ws.Range("R4").Value = 9.429637844282002
$ws.Range("S4").Value = 0.02196066496713707
$ws.Range("T4").Value = 0.02196066496713707

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lta"
$ws.Range("C5").Value = "Tnfrsf1b"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4993536666666667
$ws.Range("H5").Value = 1.498061
$ws.Range("I5").Value = 0.9584860631692095
$ws.Range("J5").Value = 0.9584860631692095
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 66.20156266666667
$ws.Range("N5").Value = 198.604688
$ws.Range("O5").Value = 0.7229090539378943
$ws.Range("P5").Value = 0.7229090539378942
$ws.Range("Q5").Value = 33.05799305666311
$ws.Range("R5").Value = 297.521937509968
$ws.Range("S5").Value = 0.69289825313831
$ws.Range("T5").Value = 0.6928982531383099

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Lta"
$ws.Range("C6").Value = "Tnfrsf1b"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.021628
$ws.Range("H6").Value = 0.064884
$ws.Range("I6").Value = 0.04151393683079058
$ws.Range("J6").Value = 0.04151393683079059
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 16.05260533333333
$ws.Range("N6").Value = 48.157816
$ws.Range("O6").Value = 0.1752915379534001
$ws.Range("P6").Value = 0.1752915379534001
$ws.Range("Q6").Value = 0.3471857481493333
$ws.Range("R6").Value = 3.124671733344
$ws.Range("S6").Value = 0.007277041833569581
$ws.Range("T6").Value = 0.007277041833569582

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Lta"
$ws.Range("C7").Value = "Tnfrsf1b"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.021628
$ws.Range("H7").Value = 0.064884
$ws.Range("I7").Value = 0.04151393683079058
$ws.Range("J7").Value = 0.04151393683079059
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.224257666666666
$ws.Range("N7").Value = 21.672773
$ws.Range("O7").Value = 0.07888758308485012
$ws.Range("P7").Value = 0.07888758308485012
$ws.Range("Q7").Value = 0.1562462448146666
$ws.Range("R7").Value = 1.406216203332
$ws.Range("S7").Value = 0.003274934140918212
$ws.Range("T7").Value = 0.003274934140918212

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Lta"
$ws.Range("C8").Value = "Tnfrsf1b"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.021628
$ws.Range("H8").Value = 0.064884
$ws.Range("I8").Value = 0.04151393683079058
$ws.Range("J8").Value = 0.04151393683079059
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.098187333333334
$ws.Range("N8").Value = 6.294562000000001
$ws.Range("O8").Value = 0.02291182502385553
$ws.Range("P8").Value = 0.02291182502385553
$ws.Range("Q8").Value = 0.04537959564533334
$ws.Range("R8").Value = 0.408416360808
$ws.Range("S8").Value = 0.0009511600567184654
$ws.Range("T8").Value = 0.0009511600567184656

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Lta"
$ws.Range("C9").Value = "Tnfrsf1b"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.021628
$ws.Range("H9").Value = 0.064884
$ws.Range("I9").Value = 0.04151393683079058
$ws.Range("J9").Value = 0.04151393683079059
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 66.20156266666667
$ws.Range("N9").Value = 198.604688
$ws.Range("O9").Value = 0.7229090539378943
$ws.Range("P9").Value = 0.7229090539378942
$ws.Range("Q9").Value = 1.431807397354667
$ws.Range("R9").Value = 12.886266576192
$ws.Range("S9").Value = 0.03001080079958433
$ws.Range("T9").Value = 0.03001080079958433
